# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scrape update).
# Source cells are plain text (coinranking.com export), so numeric-looking prices
# are written with a leading apostrophe to force Excel to keep them as text
# (otherwise Excel auto-converts strings like "210.92" into numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '25.913.50'
$ws.Range('E2').Value = '  +0.04%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.590.65'
$ws.Range('E3').Value = '  -0.77%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.43%  '

# Row 5: BNB
$ws.Range('D5').Value = '''210.92'
$ws.Range('E5').Value = '  +0.43%  '

# Row 7: XRP
$ws.Range('D7').Value = '''0.477'
$ws.Range('E7').Value = '  -0.98%  '

# Row 8: Cardano
$ws.Range('D8').Value = '''0.249'
$ws.Range('E8').Value = '  +1.62%  '

# Row 9: Dogecoin
$ws.Range('E9').Value = '  +0.35%  '

# Row 10: Solana
$ws.Range('D10').Value = '''18.30'
$ws.Range('E10').Value = '  +2.52%  '

# Row 11: TRON
$ws.Range('D11').Value = '''0.0790'
$ws.Range('E11').Value = '  +0.23%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.811.80'
$ws.Range('E12').Value = '  -0.70%  '

# Row 13: WrappedEther
$ws.Range('D13').Value = '1.589.87'
$ws.Range('E13').Value = '  -0.79%  '

# Row 14: Polkadot
$ws.Range('D14').Value = '''4.04'
$ws.Range('E14').Value = '  -0.32%  '

# Row 15: Polygon
$ws.Range('D15').Value = '''0.507'
$ws.Range('E15').Value = '  -0.64%  '

# Row 16: WrappedBTC
$ws.Range('D16').Value = '25.911.32'
$ws.Range('E16').Value = '  +0.13%  '

# Row 17: ShibaInu
$ws.Range('D17').Value = '0.0₃0725'
$ws.Range('E17').Value = '  +0.26%  '

# Row 18: Litecoin
$ws.Range('D18').Value = '''60.19'
$ws.Range('E18').Value = '  -1.47%  '

# Row 19: Dai
$ws.Range('E19').Value = '  -0.39%  '

# Row 20: BitcoinCash
$ws.Range('D20').Value = '''195.23'
$ws.Range('E20').Value = '  +3.14%  '

# Row 21: Uniswap
$ws.Range('D21').Value = '''4.21'
$ws.Range('E21').Value = '  +1.03%  '

# Row 22: Avalanche
$ws.Range('E22').Value = '  +0.68%  '

# Row 23: Chainlink
$ws.Range('D23').Value = '''5.97'
$ws.Range('E23').Value = '  +0.68%  '

# Row 24: Stellar
$ws.Range('E24').Value = '  +2.08%  '

# Row 25: Monero
$ws.Range('D25').Value = '''141.72'
$ws.Range('E25').Value = '  -0.15%  '

# Row 26: BinanceUSD
$ws.Range('E26').Value = '  -0.44%  '

# Row 27: Toncoin
$ws.Range('E27').Value = '  +0.42%  '

# Row 28: EthereumClassic
$ws.Range('E28').Value = '  +1.17%  '

# Row 29: Cosmos
$ws.Range('E29').Value = '  -0.62%  '

# Row 30: PancakeSwap
$ws.Range('E30').Value = '  -3.38%  '

# Row 31: Hedera
$ws.Range('D31').Value = '''0.0474'
$ws.Range('E31').Value = '  +0.77%  '

# Row 32: Filecoin
$ws.Range('E32').Value = '  +1.93%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range('E33').Value = '  +0.26%  '

# Row 34: LidoDAOToken
$ws.Range('D34').Value = '''1.51'
$ws.Range('E34').Value = '  +2.85%  '

# Row 35: HuobiToken
$ws.Range('D35').Value = '''2.33'
$ws.Range('E35').Value = '  -2.84%  '

# Row 36: Maker
$ws.Range('D36').Value = '1.103.55'
$ws.Range('E36').Value = '  -0.61%  '

# Row 37: MXToken
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = '''2.35'
$ws.Range('E37').Value = '  -1.28%  '

# Row 38: PaxDollar
$ws.Range('B38').Value = 'PaxDollar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D38').Value = '''1.00'
$ws.Range('E38').Value = '  -0.48%  '

# Row 39: VeChain
$ws.Range('E39').Value = '  +0.38%  '

# Row 40: ImmutableX
$ws.Range('E40').Value = '  +1.46%  '

# Row 41: ARBITRUM
$ws.Range('E41').Value = '  -2.19%  '

# Row 42: TrustWalletToken
$ws.Range('D42').Value = '''0.799'
$ws.Range('E42').Value = '  +6.94%  '

# Row 43: FraxShare
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''5.13'
$ws.Range('E43').Value = '  +1.44%  '

# Row 44: Quant
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '''93.06'
$ws.Range('E44').Value = '  -2.61%  '

# Row 45: RocketPoolETH
$ws.Range('D45').Value = '1.725.44'
$ws.Range('E45').Value = '  -0.71%  '

# Row 46: BabyDogeCoin
$ws.Range('E46').Value = '  -2.44%  '

# Row 47: RenderToken
$ws.Range('D47').Value = '''1.54'
$ws.Range('E47').Value = '  +5.12%  '

# Row 48: Aave
$ws.Range('D48').Value = '''53.35'
$ws.Range('E48').Value = '  +0.28%  '

# Row 49: Cronos
$ws.Range('E49').Value = '  -0.61%  '

# Row 50: Mantle
$ws.Range('E50').Value = '  -0.70%  '

# Row 51: USDD
$ws.Range('E51').Value = '  -0.42%  '
